$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text (matching source data),
# not auto-converted to numbers by Excel, by pre-formatting those cells as Text.

$ws.Range("D2").Value = "28.907.43"
$ws.Range("E2").Value = "  +8.12%  "

$ws.Range("D3").Value = "1.812.62"
$ws.Range("E3").Value = "  +5.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.03"
$ws.Range("E5").Value = "  +3.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4950"
$ws.Range("E7").Value = "  +2.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2785"
$ws.Range("E8").Value = "  +7.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06426"
$ws.Range("E9").Value = "  +3.97%  "

$ws.Range("D10").Value = "1.807.82"
$ws.Range("E10").Value = "  +4.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.82"
$ws.Range("E11").Value = "  +5.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07076"
$ws.Range("E12").Value = "  +3.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6473"
$ws.Range("E13").Value = "  +7.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "83.98"
$ws.Range("E14").Value = "  +8.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.692"
$ws.Range("E15").Value = "  +5.25%  "

$ws.Range("D16").Value = "28.926.55"
$ws.Range("E16").Value = "  +8.93%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9991"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007361"
$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.29"
$ws.Range("E20").Value = "  +8.07%  "

$ws.Range("D21").Value = "2.040.37"
$ws.Range("E21").Value = "  +5.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.586"
$ws.Range("E22").Value = "  +4.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.838"
$ws.Range("E23").Value = "  +3.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.364"
$ws.Range("E24").Value = "  +6.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.74"
$ws.Range("E25").Value = "  +2.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "129.51"
$ws.Range("E26").Value = "  +21.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.42"
$ws.Range("E27").Value = "  +7.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.896"
$ws.Range("E28").Value = "  +6.82%  "

$ws.Range("E29").Value = "  +3.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.149"
$ws.Range("E30").Value = "  +3.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08361"
$ws.Range("E31").Value = "  +5.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.818"
$ws.Range("E32").Value = "  +4.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04960"
$ws.Range("E33").Value = "  +10.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.107"
$ws.Range("E34").Value = "  +10.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6739"
$ws.Range("E35").Value = "  +9.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.673"
$ws.Range("E36").Value = "  +3.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.325"
$ws.Range("E37").Value = "  +15.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.752"
$ws.Range("E38").Value = "  +12.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9532"
$ws.Range("E39").Value = "  +2.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.168"
$ws.Range("E40").Value = "  +9.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01594"
$ws.Range("E41").Value = "  +6.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9990"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4104"
$ws.Range("E43").Value = "  +7.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.39"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.168"
$ws.Range("E45").Value = "  +5.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1223"
$ws.Range("E46").Value = "  +6.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05521"
$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "31.78"
$ws.Range("E48").Value = "  +5.79%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.130"
$ws.Range("E49").Value = "  +2.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3640"
$ws.Range("E50").Value = "  +8.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.311"
$ws.Range("E51").Value = "  +5.49%  "
